$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new price-list date.
$ws.Name = "15062018"

# --- Row 16: new part "Радиатор АКПП, Hayden OC1405" ---
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "Радиатор АКПП, Hayden OC1405"
$ws.Cells.Item(16, 4).Value = "OC1405"
$ws.Cells.Item(16, 5).Value = "Hayden"
$ws.Cells.Item(16, 6).Value = "Hayden|1405&Hayden|1404&Hayden|1403|Hayden|1402&Hayden|1401"
$ws.Cells.Item(16, 7).Value = "Универсальный радиатор АКПП Hayden OC1405. В комплекте пластиковые стяжки для монтажа, шланг с хомутами.  Данная деталь в наличии. Оплата товара за наличный расчет."

# H16 / J16 ("is_new" / "is_available") reuse the existing text "True" shared
# string (not a boolean) - copy it from an existing cell so the value stays a
# text cell instead of Excel's auto-boolean coercion.
$ws.Cells.Item(2, 8).Copy()
$ws.Cells.Item(16, 8).PasteSpecial(-4104)
$excel.CutCopyMode = $false

$ws.Cells.Item(16, 9).Value = 3800

$ws.Cells.Item(2, 10).Copy()
$ws.Cells.Item(16, 10).PasteSpecial(-4104)
$excel.CutCopyMode = $false

$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 1

$ws.Cells.Item(16, 14).Value = "https://www.nixparts.com/assets/pictures/Hayden/1405_01.JPG,https://www.nixparts.com/assets/pictures/Hayden/1405_02.JPG,https://www.nixparts.com/assets/pictures/Hayden/1405_03.JPG"
# Column N cells in this sheet carry no explicit cell style (unlike every
# other column) - reset the format to match after the value write.
$ws.Cells.Item(2, 14).Copy()
$ws.Cells.Item(16, 14).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(16, 15).Value = "Hummer H2, Jeep Grand Cherokee, Cadillac Escalade, Ram 1500, Suzuki Jimny, Chevrolet Tahoe,  Volvo XC90, Infiniti FX, Volvo S80, Volvo S60, Mitsubishi Pajero Sport, Saab 9-5"
$ws.Cells.Item(16, 16).Value = 1

# --- Row 17: new part "Радиатор АКПП, Hayden 1405" ---
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "Радиатор АКПП, Hayden 1405"
$ws.Cells.Item(17, 4).Value = 1405
$ws.Cells.Item(17, 5).Value = "Hayden"
$ws.Cells.Item(17, 6).Value = "Hayden|1405&Hayden|1404&Hayden|1403|Hayden|1402&Hayden|1401"
$ws.Cells.Item(17, 7).Value = "Универсальный радиатор АКПП Hayden 1405. В комплекте пластиковые стяжки для монтажа, шланг с хомутами.  Данная деталь в наличии. Оплата товара за наличный расчет."

$ws.Cells.Item(2, 8).Copy()
$ws.Cells.Item(17, 8).PasteSpecial(-4104)
$excel.CutCopyMode = $false

$ws.Cells.Item(17, 9).Value = 3800

$ws.Cells.Item(2, 10).Copy()
$ws.Cells.Item(17, 10).PasteSpecial(-4104)
$excel.CutCopyMode = $false

$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 1

$ws.Cells.Item(17, 14).Value = "https://www.nixparts.com/assets/pictures/Hayden/1405_01.JPG,https://www.nixparts.com/assets/pictures/Hayden/1405_02.JPG,https://www.nixparts.com/assets/pictures/Hayden/1405_03.JPG"
$ws.Cells.Item(2, 14).Copy()
$ws.Cells.Item(17, 14).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(17, 15).Value = "Hummer H2, Jeep Grand Cherokee, Cadillac Escalade, Ram 1500, Suzuki Jimny, Chevrolet Tahoe,  Volvo XC90, Infiniti FX, Volvo S80, Volvo S60, Mitsubishi Pajero Sport, Saab 9-5"
$ws.Cells.Item(17, 16).Value = 1

# Match the author's final selection position.
[void]$ws.Range("G18").Select()
